# The "network" column (A) is not needed where it is; move it so it
# becomes the new "network_code" column, located right before
# "station_code" (which was column G and stays column G).
#
# This shifts source_id/source_lat/source_lon/source_depth_m/
# source_origin_time left by one (B:F -> A:E) and re-inserts the old
# column A's data as the new column F, renaming its header from
# "network" to "network_code".

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Columns("A").Cut()
$ws.Columns("G").Insert()

$ws.Range("F1").Value = "network_code"
